$d = $word.ActiveDocument

# Locate the paragraph that ends the "retrospectiva" note - the anchor
# after which the new end-of-report content must be inserted.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Fizemos a retrospectiva, nela foi pontuada*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Anchor paragraph not found"
}

$dash = [char]0x2013

# Create the four new (still empty/plain) paragraphs first, in document
# order, so none of them inherit bold formatting from a sibling - only
# then fill in text/bold on the specific ones that need it.
$target.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item($target.Index + 1)

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($p1.Index + 1)

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($p2.Index + 1)

$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Item($p3.Index + 1)

# 2) Bold heading paragraph: "10/04 – "
$p2.Range.Text = "10/04 $dash "
$p2.Range.Font.Bold = $true
$p2.Range.Font.BoldBi = $true

# 3) Body paragraph describing the planning meeting.
$p3.Range.Text = "Fizemos a planning do projeto e decidimos manter as mesmas squads até a próxima semana para conseguirmos acabar todas as dívidas técnicas mantidas da última sprint."

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
